$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.031.96"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "2.484.54"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.519"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("E12").Value = "  -3.21%  "
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.91%  "
$ws.Range("D15").Value = "2.873.60"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").Value = "2.488.67"
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("E17").Value = "  -3.87%  "
$ws.Range("D18").Value = "47.900.01"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = "  -4.09%  "
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("D22").Value = "0.0₃0929"
$ws.Range("E22").Value = "  -2.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "271.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("E25").Value = "  -3.35%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.63%  "
$ws.Range("E30").Value = "  -4.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.72%  "
$ws.Range("E36").Value = "  -3.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.37%  "
$ws.Range("E38").Value = "  -4.28%  "
$ws.Range("E39").Value = "  -5.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "122.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.68%  "
$ws.Range("E41").Value = "  -1.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").Value = "2.000.38"
$ws.Range("E45").Value = "  -1.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("E48").Value = "  -2.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.07%  "
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.39%  "
